# Update crypto price/volume data per Sun Jun 30 17:56:57 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.677.41"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "3.411.26"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'576.78"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "'143.92"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "'7.64"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "3.993.85"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "'27.99"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "3.418.40"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "61.736.27"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").Value = "'9.18"
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").Value = "'387.48"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "'74.24"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'0.0000115"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").Value = "'0.187"
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "'8.00"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'1.40"
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'23.44"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'5.19"
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'6.95"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'168.26"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "3.445.45"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "'28.35"
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("D40").Value = "'0.0757"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").Value = "2.495.05"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("D46").Value = "'22.80"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'2.10"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("E51").Value = "  -0.51%  "
